$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1423.0769
$ws.Range("I18").Value = 1333.3334
$ws.Range("K18").Value = 1333.3334
$ws.Range("M18").Value = -1049.3334
$ws.Range("H41").Value = 2190.0625
$ws.Range("I41").Value = 3098.7144
$ws.Range("J41").Value = 1483.3334
$ws.Range("K41").Value = 3098.7144
$ws.Range("L41").Value = 1483.3334
$ws.Range("M41").Value = -2658.7144
$ws.Range("N41").Value = -2363.3334
$ws.Range("H88").Value = 1374025.1
$ws.Range("I88").Value = 1917.6666
$ws.Range("J88").Value = 2060078.9
$ws.Range("K88").Value = 1917.6666
$ws.Range("L88").Value = 2060078.9
$ws.Range("M88").Value = -1511.6666
$ws.Range("N88").Value = -2060890.9
$ws.Range("H91").Value = 1374025.1
$ws.Range("I91").Value = 1917.6666
$ws.Range("J91").Value = 2060078.9
$ws.Range("K91").Value = 1917.6666
$ws.Range("L91").Value = 2060078.9
$ws.Range("M91").Value = -513.6666
$ws.Range("N91").Value = -2062886.9
$ws.Range("H137").Value = 1409.4889
$ws.Range("I137").Value = 1080.1111
$ws.Range("J137").Value = 1629.0741
$ws.Range("K137").Value = 3240.3333
$ws.Range("L137").Value = 4887.2223
$ws.Range("M137").Value = -690.3333000000002
$ws.Range("N137").Value = -9987.222300000001
$ws.Range("H138").Value = 622633.9
$ws.Range("I138").Value = 998.41174
$ws.Range("J138").Value = 822026.4
$ws.Range("K138").Value = 2995.23522
$ws.Range("L138").Value = 2466079.2
$ws.Range("M138").Value = 2144.76478
$ws.Range("N138").Value = -2476359.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4479.25
$ws.Range("I32").Value = 4536.9067
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 4536.9067
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -4249.9067
$ws.Range("N32").Value = -2574
$ws.Range("H74").Value = 2323.4443
$ws.Range("I74").Value = 1302.2
$ws.Range("K74").Value = 1302.2
$ws.Range("M74").Value = -428.2
$ws.Range("H77").Value = 2323.4443
$ws.Range("I77").Value = 1302.2
$ws.Range("K77").Value = 6511
$ws.Range("M77").Value = -2143
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 3240.2222
$ws.Range("I132").Value = 3284.25
$ws.Range("J132").Value = 3176.182
$ws.Range("K132").Value = 9852.75
$ws.Range("L132").Value = 9528.545999999998
$ws.Range("M132").Value = -7322.75
$ws.Range("N132").Value = -14588.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1430.0333
$ws.Range("I31").Value = 1414.7693
$ws.Range("J31").Value = 1441.7059
$ws.Range("K31").Value = 1414.7693
$ws.Range("L31").Value = 1441.7059
$ws.Range("M31").Value = -1119.7693
$ws.Range("N31").Value = -2031.7059
$ws.Range("H32").Value = 4400
$ws.Range("I32").Value = 1600
$ws.Range("K32").Value = 1600
$ws.Range("M32").Value = -1284
$ws.Range("H34").Value = 1430.0333
$ws.Range("I34").Value = 1414.7693
$ws.Range("J34").Value = 1441.7059
$ws.Range("K34").Value = 1414.7693
$ws.Range("L34").Value = 1441.7059
$ws.Range("M34").Value = -1212.7693
$ws.Range("N34").Value = -1845.7059
$ws.Range("H58").Value = 1476.6
$ws.Range("I58").Value = 1199
$ws.Range("J58").Value = 1893
$ws.Range("K58").Value = 1199
$ws.Range("L58").Value = 1893
$ws.Range("M58").Value = -996
$ws.Range("N58").Value = -2299
$ws.Range("H74").Value = 32000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 32000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 32000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -33748
$ws.Range("H77").Value = 32000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 32000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 96000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -104736
$ws.Range("H86").Value = 3729489.8
$ws.Range("J86").Value = 26938.77
$ws.Range("L86").Value = 26938.77
$ws.Range("N86").Value = -29184.77
$ws.Range("H89").Value = 3729489.8
$ws.Range("J89").Value = 26938.77
$ws.Range("L89").Value = 134693.85
$ws.Range("N89").Value = -145925.85
$ws.Range("H114").Value = 23999.5
$ws.Range("J114").Value = 23999.5
$ws.Range("L114").Value = 23999.5
$ws.Range("N114").Value = -32677.5
$ws.Range("H132").Value = 2983.5833
$ws.Range("I132").Value = 2637.5
$ws.Range("J132").Value = 3329.6667
$ws.Range("K132").Value = 7912.5
$ws.Range("L132").Value = 9989.000100000001
$ws.Range("M132").Value = -5382.5
$ws.Range("N132").Value = -15049.0001
$ws.Range("H136").Value = 1476.6
$ws.Range("I136").Value = 1199
$ws.Range("J136").Value = 1893
$ws.Range("K136").Value = 3597
$ws.Range("L136").Value = 5679
$ws.Range("M136").Value = -1047
$ws.Range("N136").Value = -10779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 281.2857
$ws.Range("I23").Value = 199.5
$ws.Range("J23").Value = 314
$ws.Range("K23").Value = 598.5
$ws.Range("L23").Value = 942
$ws.Range("M23").Value = -363.5
$ws.Range("N23").Value = -1412
$ws.Range("H131").Value = 15154431
$ws.Range("J131").Value = 3429.2363
$ws.Range("L131").Value = 10287.7089
$ws.Range("N131").Value = -20367.7089

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4096.2666
$ws.Range("I132").Value = 4897.75
$ws.Range("J132").Value = 3180.2856
$ws.Range("K132").Value = 14693.25
$ws.Range("L132").Value = 9540.856800000001
$ws.Range("M132").Value = -12163.25
$ws.Range("N132").Value = -14600.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 299.6842
$ws.Range("I55").Value = 269.27274
$ws.Range("K55").Value = 269.27274
$ws.Range("M55").Value = -96.27274
$ws.Range("H100").Value = 1389
$ws.Range("I100").Value = 1389
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1389
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -848
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 3679.0667
$ws.Range("I132").Value = 4899.6
$ws.Range("J132").Value = 3068.8
$ws.Range("K132").Value = 14698.8
$ws.Range("L132").Value = 9206.400000000001
$ws.Range("M132").Value = -12168.8
$ws.Range("N132").Value = -14266.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H125").Value = 41366.125
$ws.Range("J125").Value = 41366.125
$ws.Range("L125").Value = 41366.125
$ws.Range("N125").Value = -51206.125
$ws.Range("H132").Value = 7533.385
$ws.Range("I132").Value = 11162.286
$ws.Range("J132").Value = 3299.6667
$ws.Range("K132").Value = 33486.858
$ws.Range("L132").Value = 9899.000100000001
$ws.Range("M132").Value = -30956.858
$ws.Range("N132").Value = -14959.0001
$ws.Range("H136").Value = 1093.7812
$ws.Range("I136").Value = 1100.4286
$ws.Range("K136").Value = 3301.2858
$ws.Range("M136").Value = -751.2857999999997
